$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row 7 with only A7 and D7 populated
$ws.Cells.Item(7, 1).Value = "table"
$ws.Cells.Item(7, 4).Value = "Objets"

# Rename header "Active" -> "State"
$ws.Range("C1").Value = "State"

# Convert the "Active" column from text true/false to numeric 1/0
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("C6").Value = 1

# Move the active selection to A8
$ws.Range("A8").Select()
